$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.850.60'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '3.545.89'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '203.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '555.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.51%  '
$ws.Range('D7').Value = '3.537.04'
$ws.Range('E7').Value = '  -2.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.610'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.81%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '65.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +17.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.664'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.145'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000272'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.95'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '4.110.89'
$ws.Range('E15').Value = '  -1.98%  '
$ws.Range('D16').Value = '3.546.50'
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').Value = '67.614.21'
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.98%  '
$ws.Range('E21').Value = '  -3.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '395.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.49%  '
$ws.Range('E23').Value = '  -4.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -10.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '693.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.28'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.91'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.33'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('E35').Value = '  -4.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '39.40'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.409'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.43%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  -1.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.84%  '
$ws.Range('D41').Value = '3.105.85'
$ws.Range('E41').Value = '  -3.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.997'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('D43').Value = '0.0₃0696'
$ws.Range('E43').Value = '  -11.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +13.99%  '
$ws.Range('E45').Value = '  -13.12%  '
$ws.Range('E46').Value = '  +6.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0403'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.79%  '
$ws.Range('E48').Value = '  -2.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '139.22'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.82%  '
$ws.Range('E51').Value = '  -4.85%  '
